# Append the 2026-01-10 profit-run row (row 47) to the bottom of the data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds dates stored as plain text (e.g. "01/09/2026" in A46), not
# real date serials. Force the new cell to Text format before writing so the
# "01/10/2026" string isn't auto-parsed into a date number, then restore the
# cell's style to Normal (matching the unstyled data cells above it).
$ws.Range("A47").NumberFormat = "@"
$ws.Range("A47").Value = "01/10/2026"
$ws.Range("A47").Style = "Normal"

$ws.Range("B47").Value = 12897.3
$ws.Range("C47").Value = 0.2098245391966422
$ws.Range("D47").Value = 0.7901754608033578
$ws.Range("E47").Value = -118.72
$ws.Range("F47").Value = -18.7
$ws.Range("G47").Value = -20500.63
$ws.Range("H47").Value = -66.8
$ws.Range("I47").Value = -396.68
$ws.Range("J47").Value = -12.78
